$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared-string text values in row 2 (TC / credential data)
$ws.Range("A2").Value = "clzCT835"
$ws.Range("C2").Value = "upteqso29"
$ws.Range("D2").Value = "gdM#P5$8"
$ws.Range("F2").Value = "tjrCDPaN"
$ws.Range("G2").Value = "msYb"

# Update numeric Candidate ID value
$ws.Range("B2").Value = 23103004
